$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("NSAA", "position", "dhc", 60),
    @("NSAA", "position", "overall", 60),
    @("NSAA", "position", "acts", 60),
    @("NSAA", "sensorMagneticField", "dhc", 60),
    @("NSAA", "sensorMagneticField", "overall", 60),
    @("NSAA", "sensorMagneticField", "acts", 60),
    @("NSAA", "jointAngle", "dhc", 60),
    @("NSAA", "jointAngle", "overall", 60),
    @("NSAA", "jointAngle", "acts", 60),
    @("NSAA", "AD", "dhc", 10),
    @("NSAA", "AD", "overall", 10),
    @("NSAA", "AD", "acts", 10),
    @("NMB", "position", "dhc", 60),
    @("NMB", "position", "overall", 60),
    @("NMB", "position", "acts", 60),
    @("NMB", "sensorMagneticField", "dhc", 60),
    @("NMB", "sensorMagneticField", "overall", 60),
    @("NMB", "sensorMagneticField", "acts", 60),
    @("NMB", "jointAngle", "dhc", 60),
    @("NMB", "jointAngle", "overall", 60),
    @("NMB", "jointAngle", "acts", 60),
    @("NMB", "AD", "dhc", 10),
    @("NMB", "AD", "overall", 10),
    @("NMB", "AD", "acts", 10),
    @("NMB", "AD", "dhc", 10),
    @("NMB", "AD", "overall", 10),
    @("NMB", "AD", "acts", 10)
)

$startRow = 1376
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rec = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
}
